# Insert a new slide "Criterios de selecao" right before the "Resultados"
# slide (which currently sits at index 11), pushing Resultados / Referencias
# down by one position.

$p = $ppt.ActivePresentation

$resultadosSlide = $p.Slides.Item(11)
$newSlide = $p.Slides.Add(12, $resultadosSlide.Layout)

# --- Title -----------------------------------------------------------
$title = $newSlide.Shapes.Item(1)
$titleRange = $title.TextFrame.TextRange
$titleRange.Text = "Critérios"
$titleRange2 = $titleRange.InsertAfter(" de ")
$titleRange2.InsertAfter("seleção")

# --- Table (replaces the "Content Placeholder" body) -----------------
$contentPh = $newSlide.Shapes.Item(2)
$contentPh.Delete()

$tableShape = $newSlide.Shapes.AddTable(12, 2, 48.72, 125.13, 828.0, 399.6)
$tableShape.Name = "Content Placeholder 6"
$table = $tableShape.Table

for ($r = 1; $r -le 12; $r++) {
    $table.Rows.Item($r).Height = 29.2
}
$table.Columns.Item(1).Width = 414.0
$table.Columns.Item(2).Width = 414.0

# Row 1 - header
$c = $table.Cell(1,1).Shape.TextFrame.TextRange
$c.Text = "Inclusão"
$c.ParagraphFormat.Alignment = 2
$c = $table.Cell(1,2).Shape.TextFrame.TextRange
$c.Text = "Exclusão"
$c.ParagraphFormat.Alignment = 2

# Row 2
$c = $table.Cell(2,1).Shape.TextFrame.TextRange
$c.Text = "(I) Os documentos devem estar disponíveis na web;"
$c.Font.Size = 12
$c = $table.Cell(2,2).Shape.TextFrame.TextRange
$c.Text = "(E) Nao utiliza o codigo fonte como entrada de dados para o treino supervisionado;"

# Row 3
$c = $table.Cell(3,1).Shape.TextFrame.TextRange
$c.Text = "(I) Estudos sobre uso de redes neurais na identificação e/ou correção de erros;"
$c.Font.Size = 12
$c = $table.Cell(3,2).Shape.TextFrame.TextRange
$c.Text = "(E) Nao corresponde a area proposta da pesquisa;"

# Row 4
$c = $table.Cell(4,1).Shape.TextFrame.TextRange
$c.Text = "(I) Estudos sobre uso de inteligência artificial no ensino de programação; "
$c.Font.Size = 12
$c = $table.Cell(4,2).Shape.TextFrame.TextRange
$c.Text = "(E) Aplicacao de inteligencia artificial para prever o desempenho do aluno;"

# Row 5
$c = $table.Cell(5,1).Shape.TextFrame.TextRange
$c.Text = "(I) Estudos sobre como os dados são representados num modelo de rede neural recorrente;"
$c.Font.Size = 12
$c = $table.Cell(5,2).Shape.TextFrame.TextRange
$c.Text = "(E) Proceedings;"

# Row 6
$c = $table.Cell(6,1).Shape.TextFrame.TextRange
$c.Text = "(I) Estudos sobre a representação de textos em redes recorrentes;"
$c.Font.Size = 12

# Row 7
$c = $table.Cell(7,1).Shape.TextFrame.TextRange
$c.Text = "(I) Publicações apresentam palavra da string de busca no seu título ou resumo;"
$c.Font.Size = 12

# Row 8
$c = $table.Cell(8,1).Shape.TextFrame.TextRange
$c.Text = "(I) Uso de sequencia de tokens para representar o codigo fonte;"
$c.Font.Size = 12

# Row 9
$c = $table.Cell(9,1).Shape.TextFrame.TextRange
$c.Text = "(I) Uso de AST para representar o codigo fonte numa rede neural;"
$c.Font.Size = 12

# Row 10
$c = $table.Cell(10,1).Shape.TextFrame.TextRange
$c.Text = "(I) Analise do codigo fonte por redes neurais convolucionais;"
$c.Font.Size = 12

# Row 11
$c = $table.Cell(11,1).Shape.TextFrame.TextRange
$c.Text = "(I) Analise do codigo fonte por redes neurais recorrentes;"
$c.Font.Size = 12

# Row 12
$c = $table.Cell(12,1).Shape.TextFrame.TextRange
$c.Text = "(I) Estudos sobre a representacao de codigos fontes;"
$c.Font.Size = 12
